$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add summary MAX formulas below the existing table (rows 28-29)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("N28").Formula = "=MAX(N20:N25)"
$ws1.Range("F29").Formula = "=MAX(D20:M25)"

# ---------------------------------------------------------------------------
# Sheet2: add averages / max / min summary block (rows 10-17)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("F10").Formula = "=AVERAGE(D6:K8)"

$ws2.Range("D11").Formula = "=AVERAGE(D6:D8)"
$ws2.Range("E11:K11").Formula = "=AVERAGE(E6:E8)"

$ws2.Range("D13").Formula = "=MAX(D6:K8)"

$ws2.Range("E14").Formula = "=MIN(D6:K8)"

$ws2.Range("B16").Value = 0.75166666666666604
$ws2.Range("B17").Value = 0.72333333333333305

# ---------------------------------------------------------------------------
# Selections / active sheet / view tweaks
# ---------------------------------------------------------------------------

# Sheet1: clear the frozen/scrolled topLeftCell and move the cell selection
[void]$ws1.Range("N28").Select()

# Sheet2 becomes the active (selected) sheet/tab, with a new selection
[void]$ws2.Range("J17").Select()
$ws2.Activate()
